$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows (5, 6, 7) that mirror the shape/formatting of the
# existing data rows (2-4): copy each existing row's formats down into the
# new row, then overwrite with the new row's values.

$ws.Range("A2:E2").Copy()
$ws.Range("A5:E5").PasteSpecial(-4122)

$ws.Range("A3:E3").Copy()
$ws.Range("A6:E6").PasteSpecial(-4122)

$ws.Range("A4:E4").Copy()
$ws.Range("A7:E7").PasteSpecial(-4122)

$ws.Range("A5").Value = 87
$ws.Range("B5").Value = "assessment_key_5"
$ws.Range("C5").Value = "BHI"
$ws.Range("D5").Value = 42667.840555555558
$ws.Range("E5").Value = 5

$ws.Range("A6").Value = 113
$ws.Range("B6").Value = "assessment_key_6"
$ws.Range("C6").Value = "BHI"
$ws.Range("D6").Value = 42682.545659722222
$ws.Range("E6").Value = 5

$ws.Range("A7").Value = 139
$ws.Range("B7").Value = "assessment_key_7"
$ws.Range("C7").Value = "BHI"
$ws.Range("D7").Value = 42712.627893518518
$ws.Range("E7").Value = 5

[void]$ws.Range("E10").Select()
